$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '56.368.89'
$ws.Cells.Item(2, 5).Value = '  -0.90%  '
$ws.Cells.Item(3, 4).Value = '3.006.84'
$ws.Cells.Item(3, 5).Value = '  +1.25%  '
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '1.00'
$ws.Cells.Item(4, 5).Value = '  -0.02%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '506.44'
$ws.Cells.Item(5, 5).Value = '  +1.20%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '138.65'
$ws.Cells.Item(6, 5).Value = '  +0.63%  '
$ws.Cells.Item(7, 5).Value = '  -0.02%  '
$ws.Cells.Item(8, 5).Value = '  +0.79%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '7.13'
$ws.Cells.Item(10, 5).Value = '  +0.36%  '
$ws.Cells.Item(11, 5).Value = '  +2.75%  '
$ws.Cells.Item(12, 4).Value = '3.508.90'
$ws.Cells.Item(12, 5).Value = '  +0.96%  '
$ws.Cells.Item(13, 5).Value = '  -2.04%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '26.25'
$ws.Cells.Item(14, 5).Value = '  +1.53%  '
$ws.Cells.Item(15, 5).Value = '  +1.92%  '
$ws.Cells.Item(16, 4).Value = '56.344.31'
$ws.Cells.Item(16, 5).Value = '  -1.09%  '
$ws.Cells.Item(17, 5).Value = '  -0.22%  '
$ws.Cells.Item(18, 4).Value = '2.998.76'
$ws.Cells.Item(18, 5).Value = '  +1.10%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '13.05'
$ws.Cells.Item(19, 5).Value = '  +3.67%  '
$ws.Cells.Item(20, 5).Value = '  +2.46%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '329.59'
$ws.Cells.Item(21, 5).Value = '  +3.14%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '1.00'
$ws.Cells.Item(22, 5).Value = '  +0.23%  '
$ws.Cells.Item(23, 5).Value = '  +2.47%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '64.88'
$ws.Cells.Item(24, 5).Value = '  +2.71%  '
$ws.Cells.Item(25, 4).Value = '3.105.60'
$ws.Cells.Item(25, 5).Value = '  +0.60%  '
$ws.Cells.Item(26, 5).Value = '  +0.12%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '0.162'
$ws.Cells.Item(27, 5).Value = '  -1.49%  '
$ws.Cells.Item(28, 4).Value = '0.0₃0912'
$ws.Cells.Item(28, 5).Value = '  +2.36%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '6.50'
$ws.Cells.Item(29, 5).Value = '  -0.38%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '7.04'
$ws.Cells.Item(30, 5).Value = '  -0.12%  '
$ws.Cells.Item(31, 5).Value = '  +1.26%  '
$ws.Cells.Item(32, 5).Value = '  +1.11%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '20.33'
$ws.Cells.Item(33, 5).Value = '  +1.06%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '152.61'
$ws.Cells.Item(34, 5).Value = '  -1.51%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '4.57'
$ws.Cells.Item(35, 5).Value = '  -0.64%  '
$ws.Cells.Item(36, 5).Value = '  +0.92%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '25.69'
$ws.Cells.Item(37, 5).Value = '  +6.35%  '
$ws.Cells.Item(38, 5).Value = '  +1.77%  '
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.0663'
$ws.Cells.Item(39, 5).Value = '  -0.58%  '
$ws.Cells.Item(40, 4).Value = '3.037.18'
$ws.Cells.Item(40, 5).Value = '  +1.15%  '
$ws.Cells.Item(41, 5).Value = '  -2.31%  '
$ws.Cells.Item(42, 5).Value = '  +0.02%  '
$ws.Cells.Item(43, 5).Value = '  +2.49%  '
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.655'
$ws.Cells.Item(44, 5).Value = '  +2.81%  '
$ws.Cells.Item(45, 4).Value = '2.179.39'
$ws.Cells.Item(45, 5).Value = '  -1.00%  '
$ws.Cells.Item(46, 5).Value = '  -1.90%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '5.93'
$ws.Cells.Item(47, 5).Value = '  +0.11%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '0.934'
$ws.Cells.Item(48, 5).Value = '  -0.99%  '
$ws.Cells.Item(49, 2).Value = 'VeChain'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.0239'
$ws.Cells.Item(49, 5).Value = '  +1.72%  '
$ws.Cells.Item(50, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '19.83'
$ws.Cells.Item(50, 5).Value = '  +3.25%  '
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.0857'
$ws.Cells.Item(51, 5).Value = '  -2.18%  '
